$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Проблема" labels (column B) ---
# Written in this specific order so the rebuilt shared-strings table lands
# in the same order as the target workbook.
$ws.Range("B11").Value = "Проблемы с авторскими правами на существующих соц. Сетях"
$ws.Range("B2").Value  = "Проблема поиска юмористических изображений"
$ws.Range("B3").Value  = "Проблема поиска юмористических видео"
$ws.Range("B10").Value = "Проблемы авторов с поиском аудитории"
$ws.Range("B4").Value  = "Плохая адаптированость существующих соц. Сетей для мемов"
$ws.Range("B5").Value  = "Трудности при обсуждении мемов"
$ws.Range("B6").Value  = "Неудобство общения в соц. Сетях с видео мемами"
$ws.Range("B7").Value  = "Сложность поиска определенного мема"
$ws.Range("B9").Value  = "Сложность модификавции мемов"
$ws.Range("B8").Value  = "Ограниченные возможности следить за новыми мемами"

# --- Update case counts (column C) ---
$ws.Range("C2").Value = 29
$ws.Range("C3").Value = 25
$ws.Range("C4").Value = 22
$ws.Range("C5").Value = 21
$ws.Range("C6").Value = 19
$ws.Range("C7").Value = 17
$ws.Range("C8").Value = 15
$ws.Range("C9").Value = 9
$ws.Range("C10").Value = 10
$ws.Range("C11").Value = 12

# --- Update cumulative-percentage column (column D) ---
# Some rows keep their formula, others become hard-coded values.
$ws.Range("D2").Value = 0.32
$ws.Range("D3").Formula = "=C3/SUM(C2:C10)+D2"
$ws.Range("D4").Formula = "=C4/SUM(C4:C12)"
$ws.Range("D5").Value = 0.22
$ws.Range("D6").Formula = "=C6/SUM(C6:C14)"
$ws.Range("D7").Formula = "=C7/SUM(C6:C14)+D6"
$ws.Range("D8").Formula = "=C8/SUM(C8:C16)"
$ws.Range("D9").Formula = "=C9/SUM(C8:C16)+D8"
$ws.Range("D10").Formula = "=C10/SUM(C10:C18)"
$ws.Range("D11").Value = 0.9

# --- Sheet view: zoom + active selection ---
[void]$ws.Range("B8").Select()
$excel.ActiveWindow.Zoom = 133
